# Regenerate the sheet the way the GitHub Actions workflow does: insert a
# bold/boxed "file 1" / "file 2" header row above the existing two columns
# of numbers (pushing everything down one row), and refresh a couple of the
# data values that changed at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force the cell to stay a text value (matches the source data, which
    # stores "2", "3", ... as strings, not numbers) and then drop back to
    # the workbook's default "Normal" style so no stray formatting sticks.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Final grid (row 1 is the new header; rows 2-6 are the old rows 1-5,
# shifted down, with B5/B6 updated to their new values).
Set-TextCell 2 1 "2"
Set-TextCell 2 2 "1"

Set-TextCell 3 1 "3"
Set-TextCell 3 2 "4"

Set-TextCell 4 1 "4"
Set-TextCell 4 2 "3"

Set-TextCell 5 1 "5"
Set-TextCell 5 2 "9"

$ws.Cells.Item(6, 1).Value = "הדס"
Set-TextCell 6 2 "5"

# --- New header row, bold + thin box border + centered/top-aligned ------
$ws.Cells.Item(1, 1).Value = "file 1"
$ws.Cells.Item(1, 2).Value = "file 2"

$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108  # xlCenter
$header.VerticalAlignment = -4160    # xlTop

Write-Host "header/data rewrite complete"
